# Update cryptocurrency price/volume figures in the "cryptos" worksheet.
# Each row (2-51) represents one coin; column D holds the Price and
# column E holds the 1h Volume change percentage. Only the cells whose
# values actually changed are touched below.
#
# Numeric-looking price strings (e.g. "592.46") are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original text-formatted Price column) instead of auto-converting them
# to numbers. Prices that already contain multiple "." separators
# (e.g. "68.147.43") are never auto-converted, so no prefix is needed
# there, and the Volume(1h) column values always contain spaces/% so
# they are safe to assign directly as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.147.43"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.542.42"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'592.46"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'173.90"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "2.541.81"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'0.167"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "'0.342"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'26.53"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "2.985.58"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "67.988.00"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'2.33"
$ws.Range("E18").Value = "  +131.22%  "
$ws.Range("D19").Value = "2.531.61"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'11.85"
$ws.Range("E20").Value = "  +3.67%  "
$ws.Range("D21").Value = "'8.03"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "'372.84"
$ws.Range("E22").Value = "  +4.80%  "
$ws.Range("D23").Value = "'4.15"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'4.58"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "'72.01"
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'1.92"
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("D28").Value = "'9.93"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("D29").Value = "2.669.64"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "0.0₃0970"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "'541.14"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").Value = "'8.42"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "'0.129"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'157.87"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").Value = "'19.27"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("D40").Value = "'18.63"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'39.41"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("D48").Value = "'148.14"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "'3.72"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Value = "'0.552"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  +1.18%  "
